$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# Shape 1 ("Title 1", id=2): rotate 180 deg, reposition/resize, autofit text
$shp1 = $s.Shapes.Item(1)
$shp1.Rotation = 180
$shp1.Left   = -31.636301040649414
$shp1.Top    = 492.32379150390625
$shp1.Width  = 667.0908203125
$shp1.Height = 216.22149658203125
$shp1.TextFrame.AutoSize = 2

# Shape 2 ("Text Placeholder 2", id=3): rotate 180 deg, reposition/resize
$shp2 = $s.Shapes.Item(2)
$shp2.Rotation = 180
$shp2.Left   = 32.72724533081055
$shp2.Top    = 441.5964660644531
$shp2.Width  = 448.9091491699219
$shp2.Height = 41.454490661621094
